# Daily scrape update — 2025-12-25 03:36:57 UTC
# Updates row data (rows 2-9) with freshly scraped opportunity listings,
# applies a highlight (yellow fill) on the newly-premium row's PREMIUM cell,
# and resizes several columns to fit the new content.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row data updates ---
$cell = $ws.Range("A2")
$cell.Formula = "=""1328491"""
$cell.Copy()
$cell.PasteSpecial(-4163)
$ws.Range("B2").Value = "https://aiesec.org/opportunity/global-talent/1328491"
$ws.Range("C2").Value = "[EXP] Project Support Officer"
$ws.Range("D2").Value = "Maastricht, Netherlands"
$ws.Range("E2").Value = "Yes"
$ws.Range("F2").Value = "156 applicants"
$ws.Range("G2").Value = "6 - 18 Months"
$ws.Range("H2").Value = "DHL Group"

$cell = $ws.Range("A3")
$cell.Formula = "=""1330854"""
$cell.Copy()
$cell.PasteSpecial(-4163)
$ws.Range("B3").Value = "https://aiesec.org/opportunity/global-talent/1330854"
$ws.Range("C3").Value = "Digital Production and e-commerce and Digital Assistance"
$ws.Range("D3").Value = "Αθήνα, Ελλάδα"
$ws.Range("F3").Value = "0 applicants"
$ws.Range("H3").Value = "Nanodomi"

$cell = $ws.Range("A4")
$cell.Formula = "=""1330851"""
$cell.Copy()
$cell.PasteSpecial(-4163)
$ws.Range("B4").Value = "https://aiesec.org/opportunity/global-talent/1330851"
$ws.Range("C4").Value = "Business Development Intern"
$ws.Range("D4").Value = "Lahore, Pakistan"
$ws.Range("F4").Value = "0 applicants"
$ws.Range("G4").Value = "3 - 6 Months"
$ws.Range("H4").Value = "Chughtaiz"

$cell = $ws.Range("A5")
$cell.Formula = "=""1330746"""
$cell.Copy()
$cell.PasteSpecial(-4163)
$ws.Range("B5").Value = "https://aiesec.org/opportunity/global-talent/1330746"
$ws.Range("C5").Value = "Software engineer"
$ws.Range("D5").Value = "Delft, Nederland"
$ws.Range("F5").Value = "8 applicants"
$ws.Range("H5").Value = "DENSsolutions B.V."

$cell = $ws.Range("A6")
$cell.Formula = "=""1330696"""
$cell.Copy()
$cell.PasteSpecial(-4163)
$ws.Range("B6").Value = "https://aiesec.org/opportunity/global-talent/1330696"
$ws.Range("C6").Value = "Producer"
$ws.Range("D6").Value = "6th of October City, Giza Governorate, Egypt"
$ws.Range("H6").Value = "Digitology"

$cell = $ws.Range("A7")
$cell.Formula = "=""1330695"""
$cell.Copy()
$cell.PasteSpecial(-4163)
$ws.Range("B7").Value = "https://aiesec.org/opportunity/global-talent/1330695"
$ws.Range("C7").Value = "Video Editor"
$ws.Range("D7").Value = "6th of October City, Giza Governorate, Egypt"
$ws.Range("F7").Value = "0 applicants"
$ws.Range("G7").Value = "9 - 12 Weeks"
$ws.Range("H7").Value = "Digitology"

$cell = $ws.Range("A8")
$cell.Formula = "=""1327775"""
$cell.Copy()
$cell.PasteSpecial(-4163)
$ws.Range("B8").Value = "https://aiesec.org/opportunity/global-talent/1327775"
$ws.Range("C8").Value = "Accelerate Romania| Programming Intern"
$ws.Range("D8").Value = "Bucharest, Romania"
$ws.Range("F8").Value = "89 applicants"
$ws.Range("G8").Value = "9 - 12 Weeks"
$ws.Range("H8").Value = "AQUAsoft"

$cell = $ws.Range("A9")
$cell.Formula = "=""1320933"""
$cell.Copy()
$cell.PasteSpecial(-4163)
$ws.Range("B9").Value = "https://aiesec.org/opportunity/global-talent/1320933"
$ws.Range("C9").Value = "Culinary Internship Chef"
$ws.Range("D9").Value = "Hong Kong"
$ws.Range("F9").Value = "40 applicants"
$ws.Range("G9").Value = "6 - 18 Months"
$ws.Range("H9").Value = "Treehouse"


# --- Highlight the PREMIUM cell for the new premium listing (row 2) ---
$ws.Range("E2").Interior.Color = 65535

# --- Column width adjustments to fit the refreshed content ---
$ws.Columns.Item(3).ColumnWidth = 58.1666666667   # C: 46 -> 59
$ws.Columns.Item(4).ColumnWidth = 46.1666666667   # D: 56 -> 47
$ws.Columns.Item(6).ColumnWidth = 16.1666666667   # F: 16 -> 17
$ws.Columns.Item(8).ColumnWidth = 20.1666666667   # H: 56 -> 21
